$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new header columns before the existing "shortlisted" column (H),
# shifting it to J, and set new header values. Copy header formatting from
# an existing header cell so the new cells match the bold/border/alignment
# style used across the header row.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "matched_skills"
$ws.Range("I1").Value = "required_skills"
$ws.Range("J1").Value = "shortlisted"

# Update row 2 data values
$ws.Range("A2").Value = "rajeshrajgor025@gmail.com"
$ws.Range("B2").Value = "UI/UX APPLICATION JOB"
$ws.Range("C2").Value = "machine learning, r"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "resumes\Resume.pdf"
$ws.Range("F2").Value = "# JD 5: UI/UX Designer`nUI/UX Designer: Require 0–3 years of experience in user research, wireframing..."
$ws.Range("G2").Value = 53.96
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = $false
